# news.xlsx update: replace homepage article list with the latest kaldata.com
# headlines (3 brand-new articles inserted near the top / middle, 2 stale
# articles dropped, and 2 more new articles appended at the bottom so the
# sheet keeps growing with every run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Overwrite the text of the 10 existing rows with the current
#        headlines, row by row (top to bottom). Existing hyperlinks
#        (rId1..rId10) stay anchored to these same cells. -----------------

$ws.Range("A1").Value  = "HWO е новият космически телескоп за търсене на извънземен живот"
$ws.Range("A2").Value  = "Nvidia Broadcast използва ИИ за автоматично регулиране на погледа"
$ws.Range("A3").Value  = "Нефтопреработвателната компания Exxonmobil е прогнозирала точно изменението на климата, докато публично е отхвърляла това"
$ws.Range("A4").Value  = "Заглушават GPS сигналите покрай българското крайбрежие – източникът за сега е неизвестен"
$ws.Range("A5").Value  = "Проучване: 6G технологията ще се предава чрез … хората"
$ws.Range("A6").Value  = "Ето как ChatGPT вместо мен реши теста по програмиране в интервюто за започване на работа"
$ws.Range("A7").Value  = "System76 Pangolin: идеалният лаптоп за програмиста с Ryzen 7 6800U, 144 Hz дисплей и Ubuntu 22.04"
$ws.Range("A8").Value  = "Intel чупи бариерата от 6GHz с процесора Core i9-13900KS с цена 699 долара"
$ws.Range("A9").Value  = "Lenovo ThinkPad Z13: Еволюция в действие"
$ws.Range("A10").Value = "НАСА работи върху хибридната мисия до Титан и още над дузина проекти в космоса"

# --- 2. Append the 2 extra rows that complete the current homepage list ----

$ws.Range("A11").Value = "Историческото първо изстрелване на ракета от британска територия претърпя провал"
$ws.Range("A12").Value = "Срокът за кандидатстване за космическия лагер Space Camp Turkey се удължава"

# Keep the worksheet style identical to the rest of the column (the
# "Hyperlink" cell style, same as A1:A10).
$ws.Range("A11").Style = "Hyperlink"
$ws.Range("A12").Style = "Hyperlink"

# --- 3. Wire up hyperlinks for the two brand-new rows -----------------

$ws.Range("A11").Hyperlinks.Add(
    $ws.Range("A11"),
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%b8%d1%81%d1%82%d0%be%d1%80%d0%b8%d1%87%d0%b5%d1%81%d0%ba%d0%be%d1%82%d0%be-%d0%bf%d1%8a%d1%80%d0%b2%d0%be-%d0%b8%d0%b7%d1%81%d1%82%d1%80%d0%b5%d0%bb%d0%b2%d0%b0%d0%bd%d0%b5-%d0%bd%d0%b0-%d1%80-403804.html",
    [Type]::Missing,
    "open this article"
)

$ws.Range("A12").Hyperlinks.Add(
    $ws.Range("A12"),
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%ba%d0%be%d1%81%d0%bc%d0%be%d1%81/%d1%81%d1%80%d0%be%d0%ba%d1%8a%d1%82-%d0%b7%d0%b0-%d0%ba%d0%b0%d0%bd%d0%b4%d0%b8%d0%b4%d0%b0%d1%82%d1%81%d1%82%d0%b2%d0%b0%d0%bd%d0%b5-%d0%b7%d0%b0-%d0%ba%d0%be%d1%81%d0%bc%d0%b8%d1%87%d0%b5%d1%81-404162.html",
    [Type]::Missing,
    "open this article"
)

# --- 4. Make sure the used range / dimension covers the new rows ----------

$ws.Range("A1:A12").Style = $ws.Range("A1:A12").Style
